$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the two extra columns (E, F)
$ws.Range("E1").Value = "Total Inflow"
$ws.Range("F1").Value = "Total Budget"

# Fill the data rows (2-89) with the constant budget figures
$ws.Range("E2:E89").Value = 2500
$ws.Range("F2:F89").Value = 1500

# Reset the view: scroll back to the top-left and move the selection
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$null = $ws.Range("K16").Select()
